# Weekly crime data update: new crime data collected for volume 31 number 45,
# covering the week of 11/4/2024 through 11/10/2024.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update report header text (volume/issue number and reporting week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  45"
$ws.Range("C9").Value = "Report Covering the Week  11/4/2024  Through  11/10/2024"

# --- Update crime-statistics table (rows 15-33) ---
# Some cells change from a numeric value to the "no data" placeholder text (shared
# strings "0" / "***.*") or vice versa; a style-only Copy from a same-styled
# template cell (row 14, untouched by this edit) is used to switch the cell's
# type/number-format before assigning the new value.

$ws.Range("C14").Copy($ws.Range("D15"))
$ws.Range("E14").Copy($ws.Range("E15"))
$ws.Range("G15").Value = 1
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 82
$ws.Range("J16").Value = 112
$ws.Range("K16").Value = -26.785714285714
$ws.Range("L16").Value = -21.904761904761
$ws.Range("M16").Value = -31.666666666666
$ws.Range("N16").Value = -89.182058047493
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = 11.111111111111
$ws.Range("I17").Value = 137
$ws.Range("J17").Value = 162
$ws.Range("K17").Value = -15.432098765432
$ws.Range("L17").Value = 8.730158730158
$ws.Range("M17").Value = 260.526315789474
$ws.Range("N17").Value = -39.647577092511
$ws.Range("C18").Value = 3
$ws.Range("I14").Copy($ws.Range("D18"))
$ws.Range("D18").Value = 2
$ws.Range("K14").Copy($ws.Range("E18"))
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 33.333333333333
$ws.Range("I18").Value = 121
$ws.Range("J18").Value = 166
$ws.Range("K18").Value = -27.108433734939
$ws.Range("L18").Value = -0.819672131147
$ws.Range("M18").Value = -21.428571428571
$ws.Range("N18").Value = -82.258064516129
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = 27.272727272727
$ws.Range("F19").Value = 43
$ws.Range("G19").Value = 47
$ws.Range("H19").Value = -8.510638297872
$ws.Range("I19").Value = 391
$ws.Range("J19").Value = 501
$ws.Range("K19").Value = -21.956087824351
$ws.Range("L19").Value = -29.549549549549
$ws.Range("M19").Value = 11.079545454545
$ws.Range("N19").Value = -9.907834101382
$ws.Range("I14").Copy($ws.Range("C20"))
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = -71.428571428571
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 18
$ws.Range("H20").Value = -77.777777777777
$ws.Range("I20").Value = 89
$ws.Range("J20").Value = 115
$ws.Range("K20").Value = -22.608695652173
$ws.Range("L20").Value = -21.238938053097
$ws.Range("M20").Value = 11.25
$ws.Range("N20").Value = -89.943502824858
$ws.Range("C21").Value = 21
$ws.Range("E21").Value = -16
$ws.Range("F21").Value = 79
$ws.Range("G21").Value = 94
$ws.Range("H21").Value = -15.957446808510
$ws.Range("I21").Value = 826
$ws.Range("J21").Value = 1067
$ws.Range("K21").Value = -22.586691658856
$ws.Range("L21").Value = -20.270270270270
$ws.Range("M21").Value = 10.280373831775
$ws.Range("N21").Value = -72.594558725945
$ws.Range("I14").Copy($ws.Range("C22"))
$ws.Range("C22").Value = 1
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -25
$ws.Range("I22").Value = 22
$ws.Range("J22").Value = 27
$ws.Range("K22").Value = -18.518518518518
$ws.Range("L22").Value = -18.518518518518
$ws.Range("M22").Value = -26.666666666666
$ws.Range("F23").Value = 4
$ws.Range("C14").Copy($ws.Range("G23"))
$ws.Range("E14").Copy($ws.Range("H23"))
$ws.Range("I23").Value = 22
$ws.Range("K23").Value = 4.761904761904
$ws.Range("L23").Value = 22.222222222222
$ws.Range("M23").Value = 69.230769230769
$ws.Range("C24").Value = 24
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = -20
$ws.Range("F24").Value = 70
$ws.Range("G24").Value = 91
$ws.Range("H24").Value = -23.076923076923
$ws.Range("I24").Value = 997
$ws.Range("J24").Value = 1251
$ws.Range("K24").Value = -20.303756994404
$ws.Range("L24").Value = -16.218487394958
$ws.Range("M24").Value = 24.937343358396
$ws.Range("C25").Value = 15
$ws.Range("D25").Value = 22
$ws.Range("E25").Value = -31.818181818181
$ws.Range("F25").Value = 45
$ws.Range("G25").Value = 73
$ws.Range("H25").Value = -38.356164383561
$ws.Range("I25").Value = 635
$ws.Range("J25").Value = 911
$ws.Range("K25").Value = -30.296377607025
$ws.Range("L25").Value = -28.248587570621
$ws.Range("C26").Value = 4
$ws.Range("E26").Value = -20
$ws.Range("G26").Value = 17
$ws.Range("H26").Value = -11.764705882352
$ws.Range("I26").Value = 215
$ws.Range("J26").Value = 260
$ws.Range("K26").Value = -17.307692307692
$ws.Range("L26").Value = -19.776119402985
$ws.Range("M26").Value = 33.540372670807
$ws.Range("C14").Copy($ws.Range("D27"))
$ws.Range("E14").Copy($ws.Range("E27"))
$ws.Range("G27").Value = 1
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = -50
$ws.Range("L28").Value = 43.137254901960
$ws.Range("C14").Copy($ws.Range("F31"))
$ws.Range("G31").Value = 2
$ws.Range("H31").Value = -100
$ws.Range("C14").Copy($ws.Range("C33"))

